$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.406.62"
$ws.Range("E2").Value = "  +3.62%  "
$ws.Range("D3").Value = "2.008.30"
$ws.Range("E3").Value = "  +7.59%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'0.7853"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +66.18%  "
$ws.Range("D6").Value = "'259.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.3581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +24.80%  "
$ws.Range("D9").Value = "'28.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +31.79%  "
$ws.Range("D10").Value = "'0.07054"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.03%  "
$ws.Range("D11").Value = "'0.8491"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +18.07%  "
$ws.Range("D12").Value = "'0.08112"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.11%  "
$ws.Range("D13").Value = "2.007.05"
$ws.Range("E13").Value = "  +7.64%  "
$ws.Range("D14").Value = "'100.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").Value = "'5.615"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.27%  "
$ws.Range("D16").Value = "'276.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "31.414.17"
$ws.Range("E17").Value = "  +3.70%  "
$ws.Range("D18").Value = "'14.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.66%  "
$ws.Range("D19").Value = "'5.925"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.30%  "
$ws.Range("D20").Value = "'0.000007917"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.17%  "
$ws.Range("D21").Value = "2.270.72"
$ws.Range("E21").Value = "  +7.76%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'7.187"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.17%  "
$ws.Range("D25").Value = "'10.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.29%  "
$ws.Range("D26").Value = "'0.1494"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +55.13%  "
$ws.Range("D27").Value = "'163.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").Value = "'20.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.16%  "
$ws.Range("D29").Value = "'2.372"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +26.52%  "
$ws.Range("D30").Value = "'1.625"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.91%  "
$ws.Range("D31").Value = "'4.627"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.85%  "
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("D33").Value = "'4.396"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("D34").Value = "'0.05232"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.61%  "
$ws.Range("D35").Value = "'1.223"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.40%  "
$ws.Range("D36").Value = "'0.7639"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.85%  "
$ws.Range("D37").Value = "'2.813"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("D38").Value = "'0.02012"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.73%  "
$ws.Range("D39").Value = "'2.954"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.33%  "
$ws.Range("D40").Value = "'80.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.61%  "
$ws.Range("D41").Value = "'6.691"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.01%  "
$ws.Range("D42").Value = "'2.175"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.41%  "
$ws.Range("D43").Value = "'0.4746"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.99%  "
$ws.Range("D44").Value = "'0.8595"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("D45").Value = "'104.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.57%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'7.741"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.43%  "
$ws.Range("D48").Value = "'9.925"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.49%  "
$ws.Range("D49").Value = "'0.4366"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.77%  "
$ws.Range("D50").Value = "'36.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.69%  "
$ws.Range("E51").Value = "  +14.83%  "